$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> list of columns to set from 0 to 1
$changes = @{
    3  = @("G","H")
    4  = @("D","E")
    5  = @("D","E")
    6  = @("D","E")
    7  = @("H")
    8  = @("H")
    9  = @("D","E")
    10 = @("H")
    11 = @("H")
    12 = @("H")
    13 = @("H")
    14 = @("H")
    15 = @("D","E")
    16 = @("H")
    17 = @("D","E")
    18 = @("H")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
